$wb = $excel.ActiveWorkbook

# The file e9a2adcd-bdc1-4a9e-bf9b-c36761b922f3.md has been handed off for
# localization, so its status/priority/handoff-timestamp fields are updated
# on the Overview sheet as well as on each per-language sheet (row 3 in
# every sheet corresponds to this file).

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 20:12:19"
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-22 20:12:14"
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-22 20:12:19"
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
